$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.841.30'
$ws.Range("E2").Value = '  +7.64%  '

$ws.Range("D3").Value = '1.742.11'
$ws.Range("E3").Value = '  +3.86%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9971'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.40%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9990'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3749'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.50%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3405'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.50%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.25'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.191'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.78%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07496'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.65%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9957'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.34%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.435'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.18%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.087'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.78%  '

$ws.Range("D16").Value = '1.735.39'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001082'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06750'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.22%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.86'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9976'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.82%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.225'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.76%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.81'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.24%  '

$ws.Range("D24").Value = '26.753.86'
$ws.Range("E24").Value = '  +7.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.456'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.480'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +24.60%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.442'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.82%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.61'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.35%  '

$ws.Range("D30").Value = '1.930.31'
$ws.Range("E30").Value = '  +3.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '132.51'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.97%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.107'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.65%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.100'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.33%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08624'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.83%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.700'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.85%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.95'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.37%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.452'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02369'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.33%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06295'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.90%  '

$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2184'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.507'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.225'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.45%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6292'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.36%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9971'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.930'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.15%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6118'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.081'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07223'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.81%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.05'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.19%  '
